$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.39993417263031
$ws.Range("B1").Value = 1.761411309242249
$ws.Range("C1").Value = 1.955539464950562
$ws.Range("D1").Value = 2.24391508102417
$ws.Range("E1").Value = 2.740868806838989
